$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Note comment added in L27 first so the shared-string table gets this
# text as the next new unique string (index 94), matching the target file.
$ws.Range("L27").Value = "Note: Run command will not run on '\' as input string"

# Rows that receive regression-test results: J = "-" (not applicable),
# K = "pass" (Fixed).
$rows = @(24,25,28,29,32,36,37,38,39,40,41,42,43,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = "-"
    $ws.Range("K$r").Value = "pass"
}

$ws.Range("J64").Select() | Out-Null
